$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new test case (iAU_TC_ID_247) replaces the old pre-request login row
$ws.Range("A2").Value = "iAU_TC_ID_247"
$ws.Range("B2").Value = '@RegressionA Prerequisite Validation of Exam Approve"'
$ws.Range("C2").Value = "passed"

# Row 3: Blueprint approval workflow, failed
$ws.Range("A3").Value = "iAU_TC_ID_247"
$ws.Range("B3").Value = "@RegressionA Validation of Blueprint  Approval Workflow "
$ws.Range("C3").Value = "failed"

# Row 4: Exam Approve, passed
$ws.Range("A4").Value = "iAU_TC_ID_250"
$ws.Range("B4").Value = "@RegressionA Validation of Exam Approve"
$ws.Range("C4").Value = "passed"

# Row 5: Exam Reject, passed
$ws.Range("A5").Value = "iAU_TC_ID_250"
$ws.Range("B5").Value = "@RegressionA Exam Reject"
$ws.Range("C5").Value = "passed"

# Row 6 no longer exists in the updated results - remove it
$ws.Rows.Item(6).Delete()
